$d = $word.ActiveDocument

# Locate the paragraph containing the intro sentence for the overfitting solutions.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Here are the few solutions for overfitting:*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq 0) {
    throw "Could not find anchor paragraph"
}

$texts = @(
    @{ Text = "Pruning"; Bold = $true },
    @{ Text = "Pre-Pruning"; Bold = $true },
    @{ Text = "In pre-pruning, it stops the tree construction bit early. It is preferred not to split a node if its goodness measure is below a threshold value. But it’s difficult to choose an appropriate stopping point."; Bold = $false },
    @{ Text = "Post-Pruning"; Bold = $true },
    @{ Text = "In post-pruning first, it goes deeper and deeper in the tree to build a complete tree. If the tree shows the overfitting problem then pruning is done as a post-pruning step. We use a cross-validation data to check the effect of our pruning. Using cross-validation data, it tests whether expanding a node will make an improvement or not. If it shows an improvement, then we can continue by expanding that node. But if it shows a reduction in accuracy then it should not be expanded i.e, the node should be converted to a leaf node."; Bold = $false }
)

# First, create all of the new (empty) paragraphs by repeatedly inserting a
# paragraph break right after the anchor paragraph. Since each new paragraph
# break is generated from the (non-bold) anchor paragraph mark, every new
# paragraph starts out with the anchor's (non-bold) formatting rather than
# inheriting formatting from a previously-inserted sibling paragraph.
for ($k = 0; $k -lt $texts.Count; $k++) {
    $d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()
}

# Now fill in the text (and bold where required) for each of the newly
# created paragraphs, in order.
for ($k = 0; $k -lt $texts.Count; $k++) {
    $item = $texts[$k]
    $paraIndex = $anchorIndex + 1 + $k
    $newPara = $d.Paragraphs($paraIndex)
    $newPara.Range.InsertBefore($item.Text)
    if ($item.Bold) {
        $newPara.Range.Font.Bold = $true
    }
}
